# costs_ae.xlsx: update the citation key used in the "ref" column.
# "PinarBilir2016economic" -> "bilir2016economic" (matches the updated
# LaTeX/bib citation key), for the "Elevated alanine transaminase" (row 4)
# and "Elevated aspartate transaminase" (row 5) adverse-event rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "bilir2016economic"
$ws.Range("G5").Value = "bilir2016economic"

$ws.Range("F10").Select()
